# Apply the edits described by the commit "Add files via upload" to the
# "Chiffres COVID-19 Valais.xlsx" workbook: a handful of daily-tracking input
# cells are corrected/filled in, and the window's frozen-pane scroll position
# and active-cell selection are updated to reflect where the author was
# working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some columns (L, M) are formatted as Text ("@") so a plain
# `.Value = 0` assignment gets stored as the literal string "0" instead of
# the number 0 (Excel mimics literal keystrokes into a text-formatted cell).
# Toggle the number format to General for the instant of the write, then put
# it back to Text so the cell keeps its original formatting/style.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    if ($fmt -eq "@") {
        $range.NumberFormat = "General"
        $range.Value = $value
        $range.NumberFormat = $fmt
    } else {
        $range.Value = $value
    }
}

# --- Row 241: corrected daily new-case count (cascades the running total in
#     column B for every subsequent row through ~505) ---
$ws.Range("C241").Value = 754

# --- Row 504: corrected ICU patient count (and its dependent total in H) ---
$ws.Range("G504").Value = 2

# --- Row 506: corrected daily new-case count ---
$ws.Range("C506").Value = 28

# --- Row 512: newly-reported new admission flag ---
$ws.Range("F512").Value = 1

# --- Row 513 ---
$ws.Range("C513").Value = 32
$ws.Range("F513").Value = 1

# --- Row 514 ---
$ws.Range("F514").Value = 1

# --- Row 515 ---
$ws.Range("C515").Value = 13
$ws.Range("F515").Value = 1

# --- Row 516 ---
$ws.Range("C516").Value = 13
$ws.Range("F516").Value = 1

# --- Row 517 ---
$ws.Range("C517").Value = 31
$ws.Range("F517").Value = 1

# --- Row 518: newly-filled-in day ---
$ws.Range("C518").Value = 22
$ws.Range("E518").Value = 1
$ws.Range("F518").Value = 1
$ws.Range("G518").Value = 2
Set-NumericValue $ws.Range("L518") 0
Set-NumericValue $ws.Range("M518") 0

# --- Row 519: newly-filled-in day ---
$ws.Range("C519").Value = 15
$ws.Range("E519").Value = 1
$ws.Range("F519").Value = 1
$ws.Range("G519").Value = 3
Set-NumericValue $ws.Range("L519") 0
Set-NumericValue $ws.Range("M519") 0

# --- Row 520: newly-filled-in day (C520 stays blank) ---
$ws.Range("E520").Value = 1
$ws.Range("F520").Value = 1
$ws.Range("G520").Value = 3
Set-NumericValue $ws.Range("L520") 0
Set-NumericValue $ws.Range("M520") 0
$ws.Range("H520").Formula = '=IF(TODAY()>A519,G520+E520,"")'

# --- Window view: scroll the frozen pane back up near the top, and move the
#     active-cell selection to where the author left off editing ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("Q515").Select()
